$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update keyword cells: remove spaces after commas (and one missing comma in B3)
$ws.Range("B3").Value = "ejuicee juice,vape liquid,vape juice,e liquid,eliquid"
$ws.Range("C3").Value = "liquid,juice,flavor,flavour,mix,pg"
$ws.Range("B7").Value = "vape pod,pod pystem,pod,pod mod "
$ws.Range("B9").Value = "vape subscription,vape box "
$ws.Range("B10").Value = "hookah,cannabis,weed,thc,cbd,marijuana"
$ws.Range("B11").Value = "vape,vaping,smoke,vapor"

# Row 3 no longer needs the extra height since the text got shorter;
# let Excel recalculate the row height back to the default.
$ws.Rows.Item(3).AutoFit()

# Update the selected cell in the active sheet view
$ws.Range("B18").Select()
